$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 800
$ws1.Range("F3").Value = 4245
$ws1.Range("F5").Value = 768

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 1

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 800
$ws4.Range("F3").Value = 4245
$ws4.Range("F5").Value = 768
$ws4.Range("F6").Value = 1
